# Commit: "update sample files. rename of text fields"
#
# The two free-text "cause of death" header columns were renamed to match
# the naming already used for the equivalent A-D columns
# (CauseOfDeathTextA/B/C/D):
#
#   Z1  (was "TextualCauseOfDeathE")     -> "CauseOfDeathTextE"
#   AD1 (was "TextualCauseOfDeathPart2") -> "CauseOfDeathTextPart2"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename AD1 before Z1 so the two renamed labels are appended to the
# shared-string table in the same order as in the saved workbook.
$ws.Range("AD1").Value = "CauseOfDeathTextPart2"
$ws.Range("Z1").Value = "CauseOfDeathTextE"

# Restore the view/selection state the author left the sheet in.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("H10").Select() | Out-Null
